# Apply updates to "Hoja1" worksheet (GANCHO J.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date in A1 (45406 -> 45436, i.e. 2024-04-24 -> 2024-05-24)
$ws.Range("A1").Value = 45436

# Update the price column (D) values for rows 29-37
$ws.Range("D29").Value = 106.327
$ws.Range("D30").Value = 113.844
$ws.Range("D31").Value = 119.215
$ws.Range("D32").Value = 125.122
$ws.Range("D33").Value = 131.032
$ws.Range("D34").Value = 138.543
$ws.Range("D35").Value = 150.362
$ws.Range("D36").Value = 169.696
$ws.Range("D37").Value = 198.694
